$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -11.9225
$ws.Range("B3").Value = 5.554900000000001
$ws.Range("C3").Value = -12.4059
$ws.Range("B4").Value = 8.739200000000004
$ws.Range("D8").Value = -8.8225
$ws.Range("C9").Value = -10.3613
$ws.Range("A11").Value = -21.9005
$ws.Range("D11").Value = -7.449599999999995
$ws.Range("A12").Value = -21.5451
$ws.Range("B14").Value = 6.392600000000001
$ws.Range("D14").Value = -7.543499999999999
$ws.Range("A15").Value = -21.84459999999998
$ws.Range("C15").Value = -13.379
$ws.Range("D15").Value = -8.390499999999994
$ws.Range("D17").Value = -8.567899999999998
$ws.Range("C19").Value = -11.86420000000001
$ws.Range("C20").Value = -11.73460000000001
$ws.Range("C25").Value = -13.1513
$ws.Range("B26").Value = 5.487199999999997
$ws.Range("D26").Value = -8.144600000000002
$ws.Range("A27").Value = -21.39669999999998
$ws.Range("C27").Value = -13.07339999999999
$ws.Range("A28").Value = -21.65049999999998
$ws.Range("C28").Value = -13.29649999999999
$ws.Range("C30").Value = -13.76929999999999
$ws.Range("A31").Value = -21.7041
$ws.Range("B31").Value = 4.999400000000003
$ws.Range("A32").Value = -21.6994
$ws.Range("C32").Value = -13.258
$ws.Range("B35").Value = 8.689900000000002
$ws.Range("A36").Value = -20.0247
$ws.Range("D36").Value = -7.635300000000002
$ws.Range("B37").Value = 9.194000000000003
$ws.Range("A38").Value = -19.3487
$ws.Range("B39").Value = 9.266400000000006
$ws.Range("B40").Value = 8.823699999999997
$ws.Range("D42").Value = -8.699299999999999
$ws.Range("C44").Value = -13.30339999999999
$ws.Range("B45").Value = 5.849099999999998
$ws.Range("A46").Value = -21.8648
$ws.Range("C47").Value = -12.272
$ws.Range("B52").Value = 5.395100000000002
$ws.Range("A54").Value = -21.38899999999998
$ws.Range("A55").Value = -22.20750000000001
$ws.Range("A56").Value = -21.7945
$ws.Range("B57").Value = 5.084699999999997
$ws.Range("C58").Value = -12.2079
$ws.Range("C62").Value = -14.5748
$ws.Range("D64").Value = -7.508999999999991
$ws.Range("A67").Value = -21.53069999999998
$ws.Range("D68").Value = -6.971999999999995
$ws.Range("A69").Value = -21.74039999999998
$ws.Range("A72").Value = -21.9711
$ws.Range("A73").Value = -20.2075
$ws.Range("C77").Value = -11.99009999999999
$ws.Range("C78").Value = -12.1266
$ws.Range("D79").Value = -5.998000000000001
$ws.Range("B81").Value = 6.3803
$ws.Range("A83").Value = -21.48999999999998
$ws.Range("B83").Value = 5.507000000000006
$ws.Range("C84").Value = -13.97179999999999
$ws.Range("A86").Value = -21.9797
$ws.Range("C89").Value = -10.7057
$ws.Range("D89").Value = -5.855700000000001
$ws.Range("A91").Value = -21.59999999999999
$ws.Range("C91").Value = -10.8169
$ws.Range("C92").Value = -11.3054
$ws.Range("A93").Value = -21.4172
$ws.Range("C96").Value = -14.058
$ws.Range("A99").Value = -20.32739999999999
$ws.Range("B100").Value = 5.228500000000001
$ws.Range("B102").Value = 8.302099999999994
$ws.Range("C102").Value = -13.9788
